# Trading update: 2026-02-17 19:57:35
# Append trade #6 (a still-OPEN MarketMaking position) as the new row 7
# on both the "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 7

    $ws.Cells.Item($row, 1).Value = 6                 # Trade #

    # Date / Time columns: force text format BEFORE assigning so Excel's
    # date auto-detection doesn't turn the literal strings into date
    # serials, then drop back to the default "Normal" style so no stray
    # number-format style gets attached to the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"       # Date
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "19:56:08"         # Time
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"     # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"             # Side
    $ws.Cells.Item($row, 6).Value = 0.48               # Entry Price
    $ws.Cells.Item($row, 7).Value = ""                 # Exit Price (still open)
    $ws.Cells.Item($row, 8).Value = "OPEN"             # Status
    $ws.Cells.Item($row, 9).Value = 0                  # P&L %
    $ws.Cells.Item($row, 10).Value = 0                 # P&L $
    $ws.Cells.Item($row, 11).Value = 99.75             # Capital After
    $ws.Cells.Item($row, 12).Value = 0                 # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                 # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6               # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = ""                # Exit Reason (still open)
    $ws.Cells.Item($row, 17).Value = 0                 # Duration (min)
}
